# Added test cases to mercedes.
#
# 1) Update existing "AMG" sheet: rows 3/4 text changes.
# 2) Add new "Shopping" sheet (after AMG) with assertion/compare-vehicles copy.
# 3) Add new "Inventory" sheet (after Shopping) with CPO/new-vehicle copy,
#    and make it the active/selected tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- AMG sheet edits ---
$ws1.Range("A3").Value = "PERFORMANCE 2/6"
$ws1.Range("A4").Value = "Mercedes-AMG Design: explore the possibilities."

# --- Shopping sheet (new) ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Shopping"
$ws2.Range("A1").Value = "Assertions"
$ws2.Range("A1").Font.Color = 0
$ws2.Range("A2").Value = "Compare Vehicles"
$ws2.Columns.Item(1).ColumnWidth = 15.7265625
$ws2.Range("A3").Select() | Out-Null

# --- Inventory sheet (new) ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Inventory"
$ws3.Range("A1").Value = "Assertions"
$ws3.Range("A1").Font.Color = 0
$ws3.Range("A2").Value = "New Vehicle Inventory"
$ws3.Range("A3").Value = "Certified Pre-Owned Inventory"
$ws3.Range("A4").Value = "Discover the Lineup"
$ws3.Range("A5").Value = "What Makes a CPO Vehicle Certified?"
$ws3.Range("A6").Value = "The Certified Pre-Owned Limited Warranty"
$ws3.Range("A7").Value = "Build Your Deal"
$ws3.Range("A8").Value = "Special Offers & More"
$ws3.Columns.Item(1).ColumnWidth = 31.90625
$ws3.Range("A9").Select() | Out-Null

# Inventory ends up the active/selected tab.
$ws3.Activate() | Out-Null
